# Enforce inform only queries
# Strip non-"inform" annotation tags (B, Ond, Sz2+, PV-loos, VC,
# Werkwoordswoordgroep, woordgroep(onderstrepen), ...) from the Word*
# columns, and correspondingly drop the matching "I" entries from each
# row's Fases (column N) list so the counts stay in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ---
$ws.Range("C7").Value = "OndWB"
$ws.Range("E7").ClearContents()
$ws.Range("N7").Value = "I,III"

# --- Row 13 ---
$ws.Range("C13").Value = "OndWB"
$ws.Range("E13").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("G13").Value = "W"
$ws.Range("H13").ClearContents()
$ws.Range("I13").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("N13").Value = "I,I,III"

# --- Row 19 ---
$ws.Range("C19").Value = "Inv,OndWBB,OndWBVC"
$ws.Range("D19").Value = "W"
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("G19").ClearContents()
$ws.Range("N19").Value = "I,III,IV,IV"

# --- Row 25 ---
$ws.Range("C25").Value = "OndWBB"
$ws.Range("E25").ClearContents()
$ws.Range("F25").ClearContents()
$ws.Range("N25").Value = "I,IV"

# --- Row 31 ---
$ws.Range("C31").Value = "OndWB"
$ws.Range("E31").ClearContents()
$ws.Range("N31").Value = "I,III"

# --- Row 37 ---
$ws.Range("C37").Value = "OndWBB"
$ws.Range("E37").ClearContents()
$ws.Range("G37").ClearContents()
$ws.Range("I37").ClearContents()
$ws.Range("N37").Value = "I,I,IV"

# --- Row 42 ---
$ws.Range("C42").ClearContents()
$ws.Range("E42").ClearContents()
$ws.Range("N42").Value = "I"

# --- Row 48 ---
$ws.Range("D48").Value = "OndWBVC"
$ws.Range("F48").ClearContents()
$ws.Range("G48").ClearContents()
$ws.Range("N48").Value = "I,IV"

# --- Row 53 ---
$ws.Range("H53").ClearContents()
$ws.Range("I53").ClearContents()
$ws.Range("N53").ClearContents()

# --- Row 54 ---
$ws.Range("C54").Value = "Inv,OndWB"
$ws.Range("E54").ClearContents()
$ws.Range("N54").Value = "I,III,III"
